$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.555.63"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.859.41"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'234.09"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'0.4705"
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.06357"
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("D10").Value = "'17.61"
$ws.Range("E10").Value = "  +7.29%  "
$ws.Range("D11").Value = "1.820.11"
$ws.Range("E11").Value = "  -2.12%  "
$ws.Range("D12").Value = "'0.07445"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "'5.129"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").Value = "'84.94"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "'0.6315"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "30.578.11"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "'243.46"
$ws.Range("E17").Value = "  +4.35%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "'12.78"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").Value = "'0.000007369"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "'5.000"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").Value = "'6.021"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'9.314"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "'164.78"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").Value = "'18.12"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").Value = "'1.896"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'0.1019"
$ws.Range("E28").Value = "  -2.02%  "
$ws.Range("D29").Value = "'1.380"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'4.056"
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("D31").Value = "'3.869"
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("D32").Value = "'0.04927"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").Value = "'1.153"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").Value = "'0.7101"
$ws.Range("E34").Value = "  -2.64%  "
$ws.Range("E35").Value = "  +0.71%  "
$ws.Range("D36").Value = "'0.01911"
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("D37").Value = "'2.689"
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("E38").Value = "  -3.13%  "
$ws.Range("D39").Value = "'1.991"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").Value = "'105.50"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").Value = "'0.4091"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").Value = "'5.560"
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("D44").Value = "'7.281"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("D45").Value = "'63.68"
$ws.Range("E45").Value = "  +3.80%  "
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").Value = "'8.585"
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("D48").Value = "'33.49"
$ws.Range("D49").Value = "'0.05546"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("D50").Value = "'1.369"
$ws.Range("E50").Value = "  -3.14%  "
$ws.Range("E51").Value = "  -0.55%  "
